# Regenerate save_data column G ("K") values: switch from old Strike# derived
# values to newly computed K values (std/mean recalculated, s_vals written).
# This updates G2:G83 in place with the freshly calculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..83 (column G), in row order.
$newK = @(
    2,1,0,2,2,1,1,1,1,1,
    1,2,1,2,1,1,1,1,1,0,
    2,2,0,0,0,2,2,0,0,0,
    1,1,1,1,1,0,2,1,0,1,
    1,0,0,0,1,1,0,2,1,1,
    0,1,3,2,1,1,0,2,1,0,
    0,1,0,2,1,2,3,2,0,0,
    2,1,1,0,0,3,1,1,3,1,
    2,1
)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
